$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: warning note about the allowed values for jenis_pegawai
$ws.Range("G1").Value = "PERHATIAN: Kolom jenis pegawai hanya boleh diisi dengan PNS, PPPK atau HONORER"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").Interior.Color = 65535

# Widen the new column so the note is readable (stored width ends up at 73)
$ws.Columns("G").ColumnWidth = 72.16666666666667

# Protect the sheet but keep the data-entry range (A2:F500) editable
$ws.Protection.AllowEditRanges.Add("input data", $ws.Range("A2:F500"))
$ws.Protect("TEMPLATE-IMPOR-PEGAWAI")

# Leave the cursor on the new note cell area
$ws.Range("E6").Select() | Out-Null
